# AE-54: Front End service for category Mappings
# Adds a new "Excel Import / Export" feature block to the top of the
# Epics sheet (pushing the existing Income/Expenses/etc. blocks down by
# three rows) and updates the active sheet/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Epics")

# Insert three new rows above the existing "Income" block (old row 6).
$ws.Rows("6:8").Insert()

# The inserted rows should pick up the same 13.95pt custom row height
# used by the blank spacer row above them (old row 5).
$ws.Rows("6:8").RowHeight = 13.95

# Populate the new feature / story cells.
$ws.Range("B6").Value = "Excel Import / Export"
$ws.Range("C7").Value = "An excel export of records"
# Row 8 stays blank, matching the usual blank spacer row between blocks.

# Make the Epics sheet the active tab, with C8 selected.
$ws.Select()
$ws.Range("C8").Select()
